# Update the cryptos table with the latest scraped prices / 1h volume
# change percentages. Column D ("Price") values are leading-apostrophe
# prefixed so Excel stores them as text (matching the workbook's original
# inlineStr cells) instead of auto-converting number-looking strings
# (e.g. "602.83") into numeric values, which would lose formatting such
# as trailing zeros ("0.0830") or thousands-style grouping dots
# ("66.120.79").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''66.120.79'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '''3.550.18'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('D5').Value = '''602.83'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('D6').Value = '''146.40'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('D7').Value = '''3.553.40'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').Value = '''0.486'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').Value = '''7.88'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '''0.133'
$ws.Range('E11').Value = '  -2.44%  '
$ws.Range('D12').Value = '''0.408'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').Value = '''4.156.58'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('D15').Value = '''29.13'
$ws.Range('E15').Value = '  -4.26%  '
$ws.Range('D16').Value = '''3.550.57'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D18').Value = '''66.152.01'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '''10.94'
$ws.Range('E19').Value = '  -5.33%  '
$ws.Range('D20').Value = '''6.23'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '''14.68'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').Value = '''418.09'
$ws.Range('E22').Value = '  -3.21%  '
$ws.Range('D23').Value = '''0.603'
$ws.Range('E23').Value = '  -1.88%  '
$ws.Range('D24').Value = '''77.72'
$ws.Range('E24').Value = '  -2.53%  '
$ws.Range('D25').Value = '''3.695.07'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = '''9.17'
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('D29').Value = '''7.89'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').Value = '''2.47'
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = '''3.550.16'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').Value = '''24.60'
$ws.Range('E34').Value = '  -3.62%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '''1.37'
$ws.Range('E36').Value = '  -6.19%  '
$ws.Range('D37').Value = '''7.63'
$ws.Range('E37').Value = '  -2.75%  '
$ws.Range('D38').Value = '''5.37'
$ws.Range('E38').Value = '  -4.84%  '
$ws.Range('D39').Value = '''1.61'
$ws.Range('E39').Value = '  -7.14%  '
$ws.Range('D40').Value = '''174.67'
$ws.Range('E40').Value = '  -0.43%  '
$ws.Range('D41').Value = '''0.0830'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').Value = '''5.10'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').Value = '''0.868'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').Value = '''45.75'
$ws.Range('E44').Value = '  -0.54%  '
$ws.Range('D45').Value = '''1.82'
$ws.Range('E45').Value = '  -6.09%  '
$ws.Range('D46').Value = '''0.999'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').Value = '''2.45'
$ws.Range('E47').Value = '  -2.16%  '
$ws.Range('D48').Value = '''23.42'
$ws.Range('E48').Value = '  -0.22%  '
$ws.Range('D49').Value = '''7.06'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.12'
$ws.Range('E50').Value = '  -6.55%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''23.72'
$ws.Range('E51').Value = '  -6.06%  '
